$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting from the existing
# header cell G1 (bold, bordered, centered style used by other headers)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding value for H2 (numeric)
$ws.Range("H2").Value = 1
